$d = $word.ActiveDocument

# --- 1. Create the new MSC_Join_A / MSC_Join_B / MSC_Join_C paragraph
#        styles, each based on the existing MSC_Join (styleId "MSCJoin")
#        style, mirroring the font overrides already used by the
#        MSC_Paragraph_B (Traditional Chinese) / MSC_Paragraph_C (Korean)
#        styles.

$sJoinA = $d.Styles.Add("MSC_Join_A", 1)
$sJoinA.BaseStyle = "MSCJoin"

$sJoinB = $d.Styles.Add("MSC_Join_B", 1)
$sJoinB.BaseStyle = "MSCJoin"
$sJoinB.Font.NameAscii = "Noto Sans CJK TC"
$sJoinB.Font.NameFarEast = "Noto Sans CJK TC"
$sJoinB.Font.Name = "Noto Sans CJK TC"
$sJoinB.Font.NameOther = "Noto Sans CJK TC"
$sJoinB.Font.NameBi = "Noto Sans CJK TC"

$sJoinC = $d.Styles.Add("MSC_Join_C", 1)
$sJoinC.BaseStyle = "MSCJoin"
$sJoinC.Font.NameAscii = "Noto Sans CJK KR"
$sJoinC.Font.NameFarEast = "Noto Sans CJK KR"
$sJoinC.Font.Name = "Noto Sans CJK KR"
$sJoinC.Font.NameOther = "Noto Sans CJK KR"
$sJoinC.Font.NameBi = "Noto Sans CJK KR"

# --- 2. Walk every paragraph in the document in order. Each table cell
#        holds a verse paragraph styled MSC_Paragraph_A/B/C followed by
#        three MSC_Join "join" paragraphs (blank, "[...]", blank) that sit
#        between the two verse ranges. Re-point those MSC_Join paragraphs
#        to the column-specific MSC_Join_A/B/C style that matches the
#        most recently seen MSC_Paragraph_A/B/C paragraph.

$paras = $d.Paragraphs
$count = $paras.Count
$cur = ""
for ($i = 1; $i -le $count; $i++) {
    $p = $paras.Item($i)
    $styleName = $p.Style.NameLocal
    if ($styleName -eq "MSC_Paragraph_A") {
        $cur = "A"
    } elseif ($styleName -eq "MSC_Paragraph_B") {
        $cur = "B"
    } elseif ($styleName -eq "MSC_Paragraph_C") {
        $cur = "C"
    } elseif ($styleName -eq "MSC_Join") {
        $p.Style = "MSC_Join_$cur"
    }
}
